$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated forecast window: dates shift forward by one day and the
# trailing row (15) rolls in a brand-new day's data.

$ws.Range("A2").Value = 45888
$ws.Range("B2").Value = 4535.36362203551
$ws.Range("C2").Value = 3893.7609393662
$ws.Range("D2").Value = 4224
$ws.Range("E2").Value = 4904.609803
$ws.Range("F2").Value = 1.62529668044529

$ws.Range("A3").Value = 45889
$ws.Range("B3").Value = 4535.36362203551
$ws.Range("C3").Value = 3937.49874211595
$ws.Range("D3").Value = 1944
$ws.Range("E3").Value = 4904.609803
$ws.Range("F3").Value = 98.4477051283516

$ws.Range("A4").Value = 45890
$ws.Range("B4").Value = 4535.36362203551
$ws.Range("C4").Value = 3966.80363917715
$ws.Range("D4").Value = 1944
$ws.Range("E4").Value = 4904.609803
$ws.Range("F4").Value = 99.6687425059014

$ws.Range("A5").Value = 45891
$ws.Range("B5").Value = 4535.36362203551
$ws.Range("C5").Value = 3969.77500140799
$ws.Range("D5").Value = 1944
$ws.Range("E5").Value = 4904.609803
$ws.Range("F5").Value = 99.7925492655201

$ws.Range("A6").Value = 45892
$ws.Range("B6").Value = 868.670076800243
$ws.Range("C6").Value = 1545.33157087306
$ws.Range("D6").Value = 1944
$ws.Range("E6").Value = 1638.789908
$ws.Range("F6").Value = 15.4771417530339

$ws.Range("A7").Value = 45893
$ws.Range("B7").Value = 773.003570643336
$ws.Range("C7").Value = 1548.8986746701
$ws.Range("D7").Value = 1944
$ws.Range("E7").Value = 1530.070577
$ws.Range("F7").Value = 15.0819033761151

$ws.Range("A8").Value = 45894
$ws.Range("B8").Value = 5476.20428821302
$ws.Range("C8").Value = 4357.3769168379
$ws.Range("D8").Value = 1944
$ws.Range("E8").Value = 5817.598163
$ws.Range("F8").Value = 114.782116317703

$ws.Range("A9").Value = 45895
$ws.Range("B9").Value = 5476.20428821302
$ws.Range("C9").Value = 4318.6910226895
$ws.Range("D9").Value = 1944
$ws.Range("E9").Value = 5817.598163
$ws.Range("F9").Value = 113.17020406152

$ws.Range("A10").Value = 45896
$ws.Range("B10").Value = 5476.20428821302
$ws.Range("C10").Value = 4272.54046207705
$ws.Range("D10").Value = 1944
$ws.Range("E10").Value = 5817.598163
$ws.Range("F10").Value = 111.247264036001

$ws.Range("A11").Value = 45897
$ws.Range("B11").Value = 5476.20428821302
$ws.Range("C11").Value = 4237.53145025912
$ws.Range("D11").Value = 1944
$ws.Range("E11").Value = 5817.598163
$ws.Range("F11").Value = 109.788555210254

$ws.Range("A12").Value = 45898
$ws.Range("B12").Value = 5476.20428821302
$ws.Range("C12").Value = 4244.99537172878
$ws.Range("D12").Value = 1944
$ws.Range("E12").Value = 5817.598163
$ws.Range("F12").Value = 110.099551938157

$ws.Range("A13").Value = 45899
$ws.Range("B13").Value = 948.243530910975
$ws.Range("C13").Value = 1595.85810714632
$ws.Range("D13").Value = 1944
$ws.Range("E13").Value = 1809.961307
$ws.Range("F13").Value = 21.3989951348059

$ws.Range("A14").Value = 45900
$ws.Range("B14").Value = 832.128236149678
$ws.Range("C14").Value = 1564.74455404578
$ws.Range("D14").Value = 1944
$ws.Range("E14").Value = 1682.358782
$ws.Range("F14").Value = 19.623962495671

$ws.Range("A15").Value = 45901
$ws.Range("B15").Value = 5776.55841866516
$ws.Range("C15").Value = 4972.39441920484
$ws.Range("D15").Value = 2952
$ws.Range("E15").Value = 6390.755553
$ws.Range("F15").Value = 109.774648064153
